# SMARTNODES.docx (Spanish) — apply the content edit described by the
# commit diff:
#   1. Rewrite the "InstantPay" bullet paragraph (now numbered with numId=2,
#      new copy text, refreshed run/paragraph formatting — Open Sans with an
#      explicit color instead of "inherit").
#   2. Rewrite the "Will have more services added later" bullet paragraph with
#      the new "SmartRewards" copy and the same refreshed formatting.
#   3. Bump the "smarthosting" bookmark id from 0 to 1 (cosmetic id only —
#      the bookmark name/position are unchanged).
#
# We use Range.InsertXML with a full single-paragraph WordProcessingML package
# so every pPr/rPr/numPr detail lands exactly as specified (Find/Replace text
# substitution alone can't touch numPr or cleanly rebuild rPr).

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, [string]$innerParagraphXml) {
    $pkg = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes" ?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$innerParagraphXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $paragraph.Range.InsertXML($pkg)
}

# --- 1) "InstantPay" bullet -> new copy, numId 1 -> 2, refreshed formatting ---
$instantPayXml = @"
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve">InstantPay (Instant Transactions): Allows for SmartCash transactions to be locked in about a second. No risk of double spending a transaction, so the receiver can trust that transaction immediately.</w:t>
  </w:r>
</w:p>
"@

$instantPayPara = $d.Paragraphs(5)
if ($instantPayPara.Range.Text -notlike "InstantPay*") {
    throw "Paragraph 5 is not the expected InstantPay bullet (got: '$($instantPayPara.Range.Text)')"
}
Set-ParagraphXml $instantPayPara $instantPayXml

# --- 2) "Will have more services added later" bullet -> SmartRewards copy ---
$smartRewardsXml = @"
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Open Sans" w:eastAsia="Times New Roman" w:hAnsi="Open Sans" w:cs="Open Sans"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve">SmartRewards: SmartRewards are calculated by the SmartNodes to allow for distribution to be handled automatically by the block rewards.</w:t>
  </w:r>
</w:p>
"@

$smartRewardsPara = $d.Paragraphs(6)
if ($smartRewardsPara.Range.Text -notlike "Will have more services*") {
    throw "Paragraph 6 is not the expected 'Will have more services' bullet (got: '$($smartRewardsPara.Range.Text)')"
}
Set-ParagraphXml $smartRewardsPara $smartRewardsXml

# --- 3) "smarthosting" bookmark: id 0 -> 1 (name/position unchanged) ---
$bmXml = @"
<w:p>
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="264" w:lineRule="atLeast"/>
    <w:jc w:val="center"/>
    <w:textAlignment w:val="baseline"/>
    <w:outlineLvl w:val="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/>
      <w:caps/>
      <w:spacing w:val="15"/>
      <w:kern w:val="36"/>
      <w:sz w:val="51"/>
      <w:szCs w:val="51"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="1" w:name="smarthosting"/>
  <w:bookmarkEnd w:id="1"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica" w:eastAsia="Times New Roman"/>
      <w:caps/>
      <w:spacing w:val="15"/>
      <w:kern w:val="36"/>
      <w:sz w:val="51"/>
      <w:szCs w:val="51"/>
    </w:rPr>
    <w:t xml:space="preserve">WHAT DO I NEED TO HOST A SMARTNODE?</w:t>
  </w:r>
</w:p>
"@

$bm = $d.Bookmarks("smarthosting")
$bmPara = $bm.Range.Paragraphs(1)
Set-ParagraphXml $bmPara $bmXml

Write-Host "Edits applied."
